# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Numeric-looking Price values are written with a leading apostrophe so Excel
# keeps storing them as literal text (matching the sheet's existing inlineStr
# cells) instead of silently reinterpreting them as numbers and dropping
# formatting such as trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.423.78"
$ws.Range("E2").Value = "  +7.35%  "
$ws.Range("D3").Value = "2.390.77"
$ws.Range("E3").Value = "  +4.90%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'113.01"
$ws.Range("E5").Value = "  +9.83%  "
$ws.Range("D6").Value = "'318.33"
$ws.Range("E6").Value = "  +2.68%  "
$ws.Range("E7").Value = "  +3.01%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +5.19%  "
$ws.Range("D10").Value = "'42.55"
$ws.Range("E10").Value = "  +10.21%  "
$ws.Range("D11").Value = "'0.0933"
$ws.Range("E11").Value = "  +3.94%  "
$ws.Range("D12").Value = "'8.70"
$ws.Range("E12").Value = "  +6.09%  "
$ws.Range("E13").Value = "  +5.78%  "
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("D15").Value = "'15.91"
$ws.Range("E15").Value = "  +5.91%  "
$ws.Range("D16").Value = "2.753.30"
$ws.Range("E16").Value = "  +4.87%  "
$ws.Range("D17").Value = "2.387.44"
$ws.Range("E17").Value = "  +4.85%  "
$ws.Range("D18").Value = "45.373.62"
$ws.Range("E18").Value = "  +6.65%  "
$ws.Range("D19").Value = "'7.68"
$ws.Range("E19").Value = "  +6.27%  "
$ws.Range("E20").Value = "  +4.29%  "
$ws.Range("D21").Value = "'13.35"
$ws.Range("E21").Value = "  +3.00%  "
$ws.Range("D22").Value = "'75.07"
$ws.Range("E23").Value = "  +5.43%  "
$ws.Range("D24").Value = "'270.02"
$ws.Range("E24").Value = "  +3.21%  "
$ws.Range("D25").Value = "'2.37"
$ws.Range("E25").Value = "  +9.38%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'11.28"
$ws.Range("E27").Value = "  +6.35%  "
$ws.Range("D28").Value = "'7.54"
$ws.Range("E28").Value = "  +9.79%  "
$ws.Range("E29").Value = "  +2.09%  "
$ws.Range("D30").Value = "'39.28"
$ws.Range("E30").Value = "  +9.84%  "
$ws.Range("D31").Value = "'22.95"
$ws.Range("E31").Value = "  +4.12%  "
$ws.Range("E32").Value = "  +11.19%  "
$ws.Range("D33").Value = "'170.47"
$ws.Range("E33").Value = "  +3.63%  "
$ws.Range("D34").Value = "'2.98"
$ws.Range("E34").Value = "  +17.16%  "
$ws.Range("D35").Value = "'0.133"
$ws.Range("E35").Value = "  +3.71%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.118"
$ws.Range("E36").Value = "  +7.68%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'4.91"
$ws.Range("E37").Value = "  +9.93%  "
$ws.Range("D38").Value = "'3.11"
$ws.Range("E38").Value = "  +17.02%  "
$ws.Range("E39").Value = "  +5.71%  "
$ws.Range("D40").Value = "'3.98"
$ws.Range("E40").Value = "  +8.75%  "
$ws.Range("D41").Value = "'1.76"
$ws.Range("E41").Value = "  +14.02%  "
$ws.Range("D42").Value = "'105.47"
$ws.Range("E42").Value = "  +7.34%  "
$ws.Range("E43").Value = "  +8.10%  "
$ws.Range("D44").Value = "'13.52"
$ws.Range("E44").Value = "  +14.27%  "
$ws.Range("D45").Value = "'71.71"
$ws.Range("E45").Value = "  +5.15%  "
$ws.Range("E46").Value = "  -0.55%  "
$ws.Range("D47").Value = "'117.71"
$ws.Range("E47").Value = "  +7.47%  "
$ws.Range("E48").Value = "  +14.18%  "
$ws.Range("D49").Value = "'1.66"
$ws.Range("E49").Value = "  +22.17%  "
$ws.Range("E50").Value = "  +9.18%  "
$ws.Range("D51").Value = "'0.223"
$ws.Range("E51").Value = "  +18.54%  "
